$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row (row 7): production_log_date, weighted_chickens_count, total_weight
$ws.Cells.Item(7, 1).Value = 45852
$ws.Cells.Item(7, 1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(7, 1).VerticalAlignment = -4108

$ws.Cells.Item(7, 2).Value = 300
$ws.Cells.Item(7, 2).VerticalAlignment = -4108

$ws.Cells.Item(7, 3).Value = 655900
$ws.Cells.Item(7, 3).VerticalAlignment = -4108

$ws.Rows.Item(7).RowHeight = 23.25

# Update the active selection to C10
$ws.Range("C10").Select()
